# Fuel Prod Imp Exp Balancing Priorities.xlsx - "updated 4.0 files and mdl"
#
# The underlying commit re-saved the workbook with a newer Excel build
# (theme name stamp, fileVersion/revisionPtr GUIDs, autofit row-height /
# column-width metrics, window geometry) alongside two real data edits and
# a cursor-position change. Only the latter are meaningful, user-driven
# changes reproducible through the Excel object model, so this script
# applies exactly those.

$wb = $excel.ActiveWorkbook

# --- About sheet: bump the "last updated" date (C1) --------------------
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Value = 45379

# --- FPIEBP sheet: re-prioritize "hard coal" (row 3) --------------------
# production / imports / exports priority columns B/C/D
# were 3/2/1 (production=3rd, imports=2nd, exports=1st priority)
# now   1/3/2 (production=1st, imports=3rd, exports=2nd priority)
$wsFPIEBP = $wb.Worksheets.Item("FPIEBP")
$wsFPIEBP.Range("B3").Value = 1
$wsFPIEBP.Range("C3").Value = 3
$wsFPIEBP.Range("D3").Value = 2

# --- FPIEBP sheet: move the active cell selection to E3 -----------------
$wsFPIEBP.Range("E3").Select() | Out-Null
